$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the plain numeric values that need to move to new cells.
#     Ordinary numbers round-trip fine through a PowerShell variable. ---
$num1 = $ws.Range("A4").Value2   # 247
$num2 = $ws.Range("B4").Value2   # 240
$num3 = $ws.Range("C4").Value2   # 244

# --- The three "5,195" / "5,234" / "5,298" labels are stored as text
#     (they contain a comma), but re-entering them through .Value would
#     get auto-coerced into numbers. Relocate them with a
#     Copy / PasteSpecial(values-only) round trip instead, which keeps
#     them as text and never touches styles.xml. Park them in a scratch
#     row far outside the table first since row/column layout below is
#     about to be rebuilt. ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("ZZ50").PasteSpecial(-4163) | Out-Null   # -4163 = xlPasteValues
$ws.Range("B2").Copy() | Out-Null
$ws.Range("ZZ51").PasteSpecial(-4163) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("ZZ52").PasteSpecial(-4163) | Out-Null

# --- Wipe the old body content (rows 2-5) so the new layout can be
#     rebuilt from a clean slate; once row 5 is completely empty it
#     drops out of sheetData entirely and the dimension shrinks back
#     down on its own. ---
$ws.Range("A2:C5").Clear()

# --- Build the new header row, re-using the existing header style
#     (bold/centered/bordered) that already lives on A1. ---
$ws.Range("A1").Copy($ws.Range("B1")) | Out-Null
$ws.Range("A1").Copy($ws.Range("C1")) | Out-Null
$ws.Range("A1").Copy($ws.Range("D1")) | Out-Null
$ws.Range("A1").Copy($ws.Range("E1")) | Out-Null

$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Germany"
$ws.Range("C1").Value = "Unnamed: 1"
$ws.Range("D1").Value = "Outside Germany"
$ws.Range("E1").Value = "Unnamed: 2"

# --- Fill in the data rows: relocated numbers go in straight away,
#     the relocated text labels are pulled back out of the scratch
#     row with another values-only paste. ---
$ws.Range("D2").Value = $num1
$ws.Range("D3").Value = $num2
$ws.Range("D4").Value = $num3

$ws.Range("ZZ50").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4163) | Out-Null
$ws.Range("ZZ51").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4163) | Out-Null
$ws.Range("ZZ52").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4163) | Out-Null

# --- Remove the scratch cells, restoring the sheet to only the final
#     A1:E4 table. ---
$ws.Range("ZZ50:ZZ52").Clear()
